$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("B1").Value = "Sample"
$ws.Range("C1").Value = "Na2O_Gl"
$ws.Range("D1").Value = "Al2O3_Gl"
$ws.Range("E1").Value = "P2O5_Gl"
$ws.Range("F1").Value = "CaO_Gl"
$ws.Range("G1").Value = "K2O_Gl"
$ws.Range("H1").Value = "TiO2_Gl"
$ws.Range("I1").Value = "SiO2_Gl"
$ws.Range("J1").Value = "MgO_Gl"
$ws.Range("K1").Value = "FeOt_Gl"
$ws.Range("L1").Value = "MnO_Gl"
$ws.Range("M1").Value = "SO2_Gl"
$ws.Range("N1").Value = "Cl_Gl"
$ws.Range("O1").Value = "Total_wt%_Gl"
$ws.Range("P1").Value = "Mg#_Gl"
$ws.Range("Q1").Value = "Na2O_STD_Gl"
$ws.Range("R1").Value = "Al2O3_STD_Gl"
$ws.Range("S1").Value = "P2O5_STD_Gl"
$ws.Range("T1").Value = "CaO_STD_Gl"
$ws.Range("U1").Value = "K2O_STD_Gl"
$ws.Range("V1").Value = "TiO2_STD_Gl"
$ws.Range("W1").Value = "SiO2_STD_Gl"
$ws.Range("X1").Value = "MgO_STD_Gl"
$ws.Range("Y1").Value = "FeOt_STD_Gl"
$ws.Range("Z1").Value = "MnO_STD_Gl"
$ws.Range("AA1").Value = "SO2_STD_Gl"
$ws.Range("AB1").Value = "Cl_STD_Gl"
$ws.Range("AC1").Value = "Total_wt%_STD_Gl"
$ws.Range("AD1").Value = "averaged?_Gl"

# --- Data row (row 2) ---
# Force "919" to stay text (Sample id), not auto-convert to a number
$ws.Range("B2").Value = "'919"
$ws.Range("C2").Value = 2.12801
$ws.Range("D2").Value = 12.511455
$ws.Range("E2").Value = 0.21128
$ws.Range("F2").Value = 10.547185
$ws.Range("G2").Value = 0.43435
$ws.Range("H2").Value = 2.41734
$ws.Range("I2").Value = 50.00389
$ws.Range("J2").Value = 8.580545
$ws.Range("K2").Value = 11.0444
$ws.Range("L2").Value = 0.37787
$ws.Range("M2").Value = 0.021205
$ws.Range("N2").Value = 0.01276
$ws.Range("O2").Value = 98.29026
$ws.Range("P2").Value = 63.31717303071684
$ws.Range("Q2").Value = 0.04709342251784689
$ws.Range("R2").Value = 0.1742315484412116
$ws.Range("S2").Value = 0.02192511497499309
$ws.Range("T2").Value = 0.1853741817244248
$ws.Range("U2").Value = 0.01052040345656425
$ws.Range("V2").Value = 0.03516656240110925
$ws.Range("W2").Value = 0.4610622624861961
$ws.Range("X2").Value = 0.4051107053427907
$ws.Range("Y2").Value = 0.148410231004021
$ws.Range("Z2").Value = 0.01706863920892479
$ws.Range("AA2").Value = 0.003451686383333353
$ws.Range("AB2").Value = 0.002361826035554317
$ws.Range("AC2").Value = 0.5226876578480366
$ws.Range("AD2").Value = "Yes"

# --- Remove now-unused trailing columns AE:AK ---
$ws.Range("AE1:AK2").Clear()
